# Tasks day 10 - actions and cookies
# Adds a new "TestAutomation" worksheet (after the existing sheets) with a
# header row: Name | Price | location

$wb = $excel.ActiveWorkbook

# Add the new worksheet at the end of the tab strip (after the last existing sheet)
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet = $wb.Worksheets.Add($null, $lastSheet)
$newSheet.Name = "TestAutomation"

# Header row for the new sheet
$newSheet.Range("A1").Value = "Name"
$newSheet.Range("B1").Value = "Price"
$newSheet.Range("C1").Value = "location"

# Keep "Sayfa2" as the active/selected sheet (unchanged from before the edit)
$wb.Worksheets.Item("Sayfa2").Activate()
